{"js": "// The document ends with a run of empty trailing paragraphs. The edit\n// collapses that run down to a single empty paragraph, and strips the\n// direct paragraph formatting on it back down to just the (redundant,\n// but explicit) \"before/after\" spacing - dropping widowControl, bidi,\n// the line/lineRule spacing attributes and the jc/alignment override.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load('items/text');\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Find the run of empty paragraphs trailing the end of the document.\nlet firstEmpty = items.length;\nwhile (firstEmpty > 0 && items[firstEmpty - 1].text === '') {\n  firstEmpty--;\n}\nconst trailingEmptyCount = items.length - firstEmpty;\n\nif (trailingEmptyCount >= 2) {\n  // Anchor on the paragraph right before the empty run (it always\n  // exists in this document - the body never starts with a blank\n  // paragraph) and insert a brand-new plain paragraph after it. A\n  // freshly split paragraph picks up the same minimal <w:pPr>/<w:rPr/>\n  // shape used throughout the rest of the document, instead of the\n  // one-off direct formatting baked onto the old last paragraph.\n  const anchor = items[firstEmpty - 1];\n  const freshParagraph = anchor.insertParagraph('', Word.InsertLocation.after);\n  await context.sync();\n\n  // Re-fetch paragraphs now that one was inserted, then drop every\n  // paragraph that belonged to the original trailing empty run (the\n  // new paragraph was inserted right before them, so they are now the\n  // last `trailingEmptyCount` items).\n  paragraphs.load('items');\n  await context.sync();\n  const refreshed = paragraphs.items;\n  const total = refreshed.length;\n  for (let i = total - trailingEmptyCount; i < total; i++) {\n    refreshed[i].delete();\n  }\n  await context.sync();\n\n  // Re-fetch once more to get a live reference to the surviving\n  // paragraph (the one we just inserted) and give it the target\n  // spacing. Office's OM only emits an explicit <w:spacing> element\n  // when the assigned value is an actual change, and before=0/after=8pt\n  // already matches the Normal style's resolved default - so the\n  // assignment is bounced off a different value first to force the\n  // engine to record (and keep) the explicit override.\n  paragraphs.load('items');\n  await context.sync();\n  const finalParagraphs = paragraphs.items;\n  const finalParagraph = finalParagraphs[finalParagraphs.length - 1];\n\n  finalParagraph.paragraphFormat.spaceBefore = 1;\n  finalParagraph.paragraphFormat.spaceAfter = 1;\n  await context.sync();\n\n  finalParagraph.paragraphFormat.spaceBefore = 0;\n  finalParagraph.paragraphFormat.spaceAfter = 8;\n  await context.sync();\n}\n", "ps1": "# The document ends with a run of empty trailing paragraphs. The edit\n# collapses that run down to a single empty paragraph, and strips the\n# direct paragraph formatting on it back down to just the (redundant,\n# but explicit) \"before/after\" spacing - dropping widowControl, bidi,\n# the line/lineRule spacing attributes and the jc/alignment override.\n\n$d = $word.ActiveDocument\n$count = $d.Paragraphs.Count\n\n# Find the last paragraph that still carries real content (its Range.Text\n# is longer than just the trailing paragraph mark); everything after it\n# is the run of blank paragraphs we need to collapse.\n$lastContentIdx = $count\nfor ($i = $count; $i -ge 1; $i--) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t.Length -gt 1) {\n        $lastContentIdx = $i\n        break\n    }\n}\n\n$oldEmptyCount = $count - $lastContentIdx\n\nif ($oldEmptyCount -ge 2) {\n    # Split a brand-new, plain paragraph in right after the last piece of\n    # real content. A freshly split paragraph picks up the same minimal\n    # <w:pPr>/<w:rPr/> shape used throughout the rest of the document,\n    # instead of the one-off direct formatting baked onto the old last\n    # paragraph.\n    $anchor = $d.Paragraphs.Item($lastContentIdx)\n    $anchor.Range.InsertParagraphAfter()\n\n    # The freshly inserted paragraph now sits right before the stale\n    # blank paragraphs, which became the new last $oldEmptyCount items -\n    # delete all of them, leaving only the fresh paragraph behind.\n    $newCount = $d.Paragraphs.Count\n    for ($i = $newCount; $i -gt ($newCount - $oldEmptyCount); $i--) {\n        $d.Paragraphs.Item($i).Range.Delete()\n    }\n\n    # Give the surviving paragraph the target spacing. The COM object\n    # model only emits an explicit <w:spacing> element when the assigned\n    # value is an actual change, and before=0/after=8pt already matches\n    # the Normal style's resolved default - so the assignment is bounced\n    # off a different value first to force it to record (and keep) the\n    # explicit override.\n    $final = $d.Paragraphs.Item($d.Paragraphs.Count)\n    $final.Format.SpaceBefore = 1\n    $final.Format.SpaceAfter = 1\n    $final.Format.SpaceBefore = 0\n    $final.Format.SpaceAfter = 8\n}\n"}
